$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 for NATLIFEINS, shifting existing rows down
$ws.Rows.Item(19).Insert()

# Re-apply the bordered/bold/centered style used by column-A ticker cells to the new row
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)

# Write the ticker name for the newly inserted row
$ws.Range("A19").Value = "NATLIFEINS"

# Update Log Return / Rsquared / Return Standard Deviation / COV for every data row (2-46)
$ws.Range("B2").Value = 0.03365286677259715
$ws.Range("C2").Value = 0.4434094891746991
$ws.Range("D2").Value = 0.8827647752603149
$ws.Range("E2").Value = 0.0441826411927045
$ws.Range("B3").Value = 0.04205036455466069
$ws.Range("C3").Value = 0.428877722450033
$ws.Range("D3").Value = 1.940962593150439
$ws.Range("E3").Value = 0.08701818380134754
$ws.Range("B4").Value = -0.004257744562832017
$ws.Range("C4").Value = 0.5880528462145973
$ws.Range("D4").Value = 0.9557605185653474
$ws.Range("E4").Value = 0.09851656850545734
$ws.Range("B5").Value = -0.01237900252246709
$ws.Range("C5").Value = 0.4532139280123447
$ws.Range("D5").Value = 0.6892420821040929
$ws.Range("E5").Value = 0.05574017365644492
$ws.Range("B6").Value = 0.02832299821265423
$ws.Range("C6").Value = 0.5461957232834214
$ws.Range("D6").Value = 0.4839040852009782
$ws.Range("E6").Value = 0.02561016006446109
$ws.Range("B7").Value = -0.03506909943976957
$ws.Range("C7").Value = 0.5967050797478487
$ws.Range("D7").Value = 0.7806155382242906
$ws.Range("E7").Value = 0.08004483894280984
$ws.Range("B8").Value = 0.01560856003690929
$ws.Range("C8").Value = 0.702372989803784
$ws.Range("D8").Value = 1.262828943615204
$ws.Range("E8").Value = 0.08327869358940596
$ws.Range("B9").Value = -0.04658961335933156
$ws.Range("C9").Value = 0.5835450676319185
$ws.Range("D9").Value = 1.018173043997194
$ws.Range("E9").Value = 0.1156376286289402
$ws.Range("B10").Value = -0.0001866199248455154
$ws.Range("C10").Value = 0.4544486333278106
$ws.Range("D10").Value = 0.8480321163425169
$ws.Range("E10").Value = 0.02130618193838865
$ws.Range("B11").Value = -0.03237355561714485
$ws.Range("C11").Value = 0.6212271026988461
$ws.Range("D11").Value = 1.635769084758878
$ws.Range("E11").Value = 0.1596697673654268
$ws.Range("B12").Value = -0.002739208666456448
$ws.Range("C12").Value = 0.6580053758317823
$ws.Range("D12").Value = 0.6882411966333436
$ws.Range("E12").Value = 0.04075121075098553
$ws.Range("B13").Value = 0.002241738847947355
$ws.Range("C13").Value = 0.4780738058840448
$ws.Range("D13").Value = 2.674900886550994
$ws.Range("E13").Value = 0.08979490452745749
$ws.Range("B14").Value = 0.05078151669859006
$ws.Range("C14").Value = 0.5411753742414527
$ws.Range("D14").Value = 1.594230917360052
$ws.Range("E14").Value = 0.07254530269460886
$ws.Range("B15").Value = -0.0377778507044137
$ws.Range("C15").Value = 0.561583663644073
$ws.Range("D15").Value = 0.7290542005722048
$ws.Range("E15").Value = 0.1068325166230784
$ws.Range("B16").Value = 0.1091447697496059
$ws.Range("C16").Value = 0.4829128939860252
$ws.Range("D16").Value = 1.620945931868722
$ws.Range("E16").Value = 0.1253699846413636
$ws.Range("B17").Value = -0.02617527650798552
$ws.Range("C17").Value = 0.4341768987414185
$ws.Range("D17").Value = 1.920308508218236
$ws.Range("E17").Value = 0.1505264882361352
$ws.Range("B18").Value = 0.0427357797455109
$ws.Range("C18").Value = 0.577585714567399
$ws.Range("D18").Value = 1.809667502758367
$ws.Range("E18").Value = 0.1544209163900706
$ws.Range("B19").Value = -0.06213329180577229
$ws.Range("C19").Value = 0.2591961132736835
$ws.Range("D19").Value = 1.560005651874994
$ws.Range("E19").Value = 0.1153467310907198
$ws.Range("B20").Value = -0.00117727708980133
$ws.Range("C20").Value = 0.5373300712835759
$ws.Range("D20").Value = 3.304403764977431
$ws.Range("E20").Value = 0.1525669338384187
$ws.Range("B21").Value = -0.03249243394983966
$ws.Range("C21").Value = 0.7190489915322813
$ws.Range("D21").Value = 0.824649283530511
$ws.Range("E21").Value = 0.08379843978419817
$ws.Range("B22").Value = 0.07561316766822147
$ws.Range("C22").Value = 0.388632313859117
$ws.Range("D22").Value = 2.769832188863254
$ws.Range("E22").Value = 0.1633843027798317
$ws.Range("B23").Value = 0.06023158032209663
$ws.Range("C23").Value = 0.465015452021884
$ws.Range("D23").Value = 2.888337713670587
$ws.Range("E23").Value = 0.1645298607042942
$ws.Range("B24").Value = 0.1393178626035699
$ws.Range("C24").Value = 0.6041289632930786
$ws.Range("D24").Value = 3.570592788968222
$ws.Range("E24").Value = 0.3226559500372305
$ws.Range("B25").Value = 0.03030515152589379
$ws.Range("C25").Value = 0.4706003796678094
$ws.Range("D25").Value = 1.209066424887155
$ws.Range("E25").Value = 0.05146265186609624
$ws.Range("B26").Value = -0.0496507203733982
$ws.Range("C26").Value = 0.5888557660704561
$ws.Range("D26").Value = 0.9563573661801972
$ws.Range("E26").Value = 0.09793877944279246
$ws.Range("B27").Value = 0.141029141430179
$ws.Range("C27").Value = 0.5794685033425042
$ws.Range("D27").Value = 2.821403228650686
$ws.Range("E27").Value = 0.1613189216394892
$ws.Range("B28").Value = -0.06003713697354095
$ws.Range("C28").Value = 0.2604495499036615
$ws.Range("D28").Value = 2.536573501112191
$ws.Range("E28").Value = 0.09888787656704598
$ws.Range("B29").Value = 0.09550446775005078
$ws.Range("C29").Value = 0.4304230546124622
$ws.Range("D29").Value = 2.295619604937464
$ws.Range("E29").Value = 0.165694049801219
$ws.Range("B30").Value = 0.4227444622695788
$ws.Range("C30").Value = 0.3311539224319301
$ws.Range("D30").Value = 4.145932991489678
$ws.Range("E30").Value = 0.2139312295939806
$ws.Range("B31").Value = 0.1922302195803633
$ws.Range("C31").Value = 0.3878717807656425
$ws.Range("D31").Value = 3.193274166938052
$ws.Range("E31").Value = 0.3877010336217573
$ws.Range("B32").Value = 0.009430889977264774
$ws.Range("C32").Value = 0.4991411996289583
$ws.Range("D32").Value = 2.363430798453433
$ws.Range("E32").Value = 0.08950252868530538
$ws.Range("B33").Value = 0.02565050196765687
$ws.Range("C33").Value = 0.4759219873944451
$ws.Range("D33").Value = 2.669258016514737
$ws.Range("E33").Value = 0.1399633027939575
$ws.Range("B34").Value = 0.2457453417264271
$ws.Range("C34").Value = 0.4836431141786909
$ws.Range("D34").Value = 3.607273073622384
$ws.Range("E34").Value = 0.3521087442698894
$ws.Range("B35").Value = 0.1163750337466325
$ws.Range("C35").Value = 0.3913409030217004
$ws.Range("D35").Value = 2.37259125782909
$ws.Range("E35").Value = 0.2460345353905922
$ws.Range("B36").Value = -0.04989808468964609
$ws.Range("C36").Value = 0.4451301447425666
$ws.Range("D36").Value = 2.57911993357996
$ws.Range("E36").Value = 0.09411748402487762
$ws.Range("B37").Value = -0.07473899551254026
$ws.Range("C37").Value = 0.5011920794541648
$ws.Range("D37").Value = 2.411133807920456
$ws.Range("E37").Value = 0.2891029487305297
$ws.Range("B38").Value = 0.3679004813935062
$ws.Range("C38").Value = 0.3432497574615242
$ws.Range("D38").Value = 3.426372234738907
$ws.Range("E38").Value = 0.2109059990473916
$ws.Range("B39").Value = -0.007409021501504159
$ws.Range("C39").Value = 0.5669670137632787
$ws.Range("D39").Value = 1.23065993185232
$ws.Range("E39").Value = 0.0670851226365345
$ws.Range("B40").Value = 0.01693803037474472
$ws.Range("C40").Value = 0.5034118854182381
$ws.Range("D40").Value = 1.996071504957383
$ws.Range("E40").Value = 0.06100048965022704
$ws.Range("B41").Value = -0.007797588426208152
$ws.Range("C41").Value = 0.5159021920008924
$ws.Range("D41").Value = 1.044795575109014
$ws.Range("E41").Value = 0.0926271976222993
$ws.Range("B42").Value = -0.006357780716271053
$ws.Range("C42").Value = 0.4775879461872391
$ws.Range("D42").Value = 1.415953229241488
$ws.Range("E42").Value = 0.0402931559843942
$ws.Range("B43").Value = -0.04553720358852453
$ws.Range("C43").Value = 0.4223239848198141
$ws.Range("D43").Value = 0.8771679215748776
$ws.Range("E43").Value = 0.0940871699621635
$ws.Range("B44").Value = -0.03885376683945015
$ws.Range("C44").Value = 0.6383915387197958
$ws.Range("D44").Value = 1.611674503947476
$ws.Range("E44").Value = 0.1388196011376104
$ws.Range("B45").Value = 0.02667274563896026
$ws.Range("C45").Value = 0.5985218973343482
$ws.Range("D45").Value = 1.000112060289442
$ws.Range("E45").Value = 0.07771600437364189
$ws.Range("B46").Value = -0.04862774278670016
$ws.Range("C46").Value = 0.4882153521599051
$ws.Range("D46").Value = 2.514337744710879
$ws.Range("E46").Value = 0.1620548983003801
